$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.702.48"
$ws.Range("E2").Value = "  -0.82%  "
$ws.Range("D3").Value = "3.089.89"
$ws.Range("E3").Value = "  -1.31%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'576.84"
$ws.Range("E5").Value = "  -0.85%  "
$ws.Range("D6").Value = "'172.74"
$ws.Range("E6").Value = "  -0.93%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "3.087.58"
$ws.Range("E8").Value = "  -1.24%  "
$ws.Range("D9").Value = "'0.514"
$ws.Range("E9").Value = "  -1.72%  "
$ws.Range("D10").Value = "'6.37"
$ws.Range("E10").Value = "  -0.58%  "
$ws.Range("E11").Value = "  -2.79%  "
$ws.Range("D12").Value = "'0.471"
$ws.Range("E12").Value = "  -2.39%  "
$ws.Range("E13").Value = "  -4.46%  "
$ws.Range("D14").Value = "'35.78"
$ws.Range("E14").Value = "  -4.50%  "
$ws.Range("E15").Value = "  -0.71%  "
$ws.Range("D16").Value = "3.603.39"
$ws.Range("E16").Value = "  -1.28%  "
$ws.Range("D17").Value = "66.601.59"
$ws.Range("E17").Value = "  -0.91%  "
$ws.Range("D18").Value = "'6.96"
$ws.Range("E18").Value = "  -2.67%  "
$ws.Range("D19").Value = "'16.84"
$ws.Range("E19").Value = "  +2.77%  "
$ws.Range("D20").Value = "3.089.95"
$ws.Range("E20").Value = "  -1.30%  "
$ws.Range("D21").Value = "'484.41"
$ws.Range("E21").Value = "  -1.95%  "
$ws.Range("D22").Value = "'7.75"
$ws.Range("E22").Value = "  -2.01%  "
$ws.Range("D23").Value = "'0.688"
$ws.Range("E23").Value = "  -3.13%  "
$ws.Range("D24").Value = "'83.32"
$ws.Range("E24").Value = "  -1.16%  "
$ws.Range("D25").Value = "'12.69"
$ws.Range("E25").Value = "  -4.90%  "
$ws.Range("E26").Value = "  -3.10%  "
$ws.Range("D27").Value = "'10.05"
$ws.Range("E27").Value = "  -3.73%  "
$ws.Range("E28").Value = "  -0.03%  "
$ws.Range("D29").Value = "'7.97"
$ws.Range("E29").Value = "  +0.34%  "
$ws.Range("D30").Value = "'2.25"
$ws.Range("E30").Value = "  -4.25%  "
$ws.Range("D32").Value = "'27.94"
$ws.Range("E32").Value = "  -2.97%  "
$ws.Range("D33").Value = "'0.112"
$ws.Range("E33").Value = "  -2.45%  "
$ws.Range("E34").Value = "  -1.51%  "
$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = "  -0.05%  "
$ws.Range("D36").Value = "'48.16"
$ws.Range("E36").Value = "  +2.52%  "
$ws.Range("D37").Value = "'5.58"
$ws.Range("E37").Value = "  -5.69%  "
$ws.Range("D38").Value = "'0.942"
$ws.Range("E38").Value = "  -3.69%  "
$ws.Range("B39").Value = "TheGraph"
$ws.Range("C39").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D39").Value = "'0.310"
$ws.Range("E39").Value = "  -0.88%  "
$ws.Range("B40").Value = "OKB"
$ws.Range("C40").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D40").Value = "'49.05"
$ws.Range("E40").Value = "  -2.21%  "
$ws.Range("D41").Value = "'0.123"
$ws.Range("E41").Value = "  -1.28%  "
$ws.Range("E42").Value = "  -5.46%  "
$ws.Range("D43").Value = "'8.28"
$ws.Range("E43").Value = "  -3.52%  "
$ws.Range("D44").Value = "'2.61"
$ws.Range("E44").Value = "  -0.16%  "
$ws.Range("D45").Value = "2.776.80"
$ws.Range("E45").Value = "  -2.13%  "
$ws.Range("D46").Value = "'0.0346"
$ws.Range("E46").Value = "  -2.44%  "
$ws.Range("D47").Value = "'368.26"
$ws.Range("E47").Value = "  -4.87%  "
$ws.Range("D48").Value = "'134.12"
$ws.Range("E48").Value = "  -1.26%  "
$ws.Range("E49").Value = "  +0.02%  "
$ws.Range("D50").Value = "'24.44"
$ws.Range("E50").Value = "  -2.46%  "
$ws.Range("E51").Value = "  -2.92%  "
